$p = $ppt.ActivePresentation

# Merge the split runs in slide 1's notes paragraph 4 into a single run.
$s1 = $p.Slides.Item(1)
$notes = $s1.NotesPage
for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
    $sh = $notes.Shapes.Item($i)
    if ($sh.Type -eq 14 -or $sh.PlaceholderFormat.Type -eq 2) {
        $tr = $sh.TextFrame.TextRange
        $tr.Paragraphs(4, 1).Text = "Les instructions du même bloc sont exécutées séquentiellement.  On retrouvera cette notion d’indentation pour les boucles for et les fonctions. "
    }
}

# Delete the trashed slides, keeping only slide 1.
for ($i = $p.Slides.Count; $i -ge 2; $i--) {
    $p.Slides.Item($i).Delete()
}
